$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column "01-nov" right before "01-oct." ---
$ws = $wb.Worksheets.Item("Prix Spot")

# Insert a brand-new column at DF; everything from the old DF ("01-oct.") through
# EJ ("31-oct.") shifts one column to the right (DF->DG, ..., EJ->EK), carrying its
# style/value along with it (matches the diff: dimension A1:EJ25 -> A1:EK25).
$ws.Range("DF1").EntireColumn.Insert()

# Header for the freshly inserted column.
$ws.Range("DF1").Value = "01-nov"

# The new day has no data yet, so every hour row (2-25) gets the placeholder "-".
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 110).Value = "-"
}

# --- Sheet "Gaz": append the new trading day 2025-10-30 ---
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Cells.Item(138, 1).NumberFormat = "@"
$wsGaz.Cells.Item(138, 1).Value = "2025-10-30"
$wsGaz.Cells.Item(138, 1).ClearFormats()
$wsGaz.Cells.Item(138, 2).Value = 29.8

# --- Sheet "CO2": append the new trading day 2025-10-30 ---
$wsCO2 = $wb.Worksheets.Item("CO2")
$wsCO2.Cells.Item(138, 1).NumberFormat = "@"
$wsCO2.Cells.Item(138, 1).Value = "2025-10-30"
$wsCO2.Cells.Item(138, 1).ClearFormats()
$wsCO2.Cells.Item(138, 2).Value = 78.36
